$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("展览") ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Rows("2:4").Delete()

# Row 2: 合肥·灵能百分百ONLY2.0
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "2024-07-27"
$ws1.Range("C2").Value = "合肥·灵能百分百ONLY2.0"
$ws1.Range("D2").Value = "铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)"
$ws1.Range("E2").Value = "2024.07.27 10:00-07.27 17:00"
$ws1.Range("F2").Value = 117
$ws1.Range("G2").Value = "不可售"
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=87497"
$ws1.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg"

# Row 3: 安徽·MAX特摄only展
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "2024-07-27"
$ws1.Range("C3").Value = "安徽·MAX特摄only展"
$ws1.Range("D3").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws1.Range("E3").Value = "2024.07.27 09:30-07.27 18:00"
$ws1.Range("F3").Value = 502
$ws1.Range("G3").Value = 50
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=83684"
$ws1.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"

# Row 4: 庐江·夏日游嘉年华
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "2024-07-27"
$ws1.Range("C4").Value = "庐江·夏日游嘉年华"
$ws1.Range("D4").Value = "白山路东150米 庐江体育馆"
$ws1.Range("E4").Value = "2024.07.27 09:00-07.28 17:00"
$ws1.Range("F4").Value = 183
$ws1.Range("G4").Value = 60
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=87569"
$ws1.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg"

# Row 5: 长丰·莓可可游戏动漫展
$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "2024-07-27"
$ws1.Range("C5").Value = "长丰·莓可可游戏动漫展"
$ws1.Range("D5").Value = "长寿路12号 长丰宾馆·梅山饭店(长寿路店)"
$ws1.Range("E5").Value = "2024.07.27 10:00-07.27 17:00"
$ws1.Range("F5").Value = 19
$ws1.Range("G5").Value = "不可售"
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87796"
$ws1.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202406/MLTfeikq1718823574810.png"

# Row 6: 合肥·咒术回战only
$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "2024-07-28"
$ws1.Range("C6").Value = "合肥·咒术回战only"
$ws1.Range("D6").Value = "清河路19号 依立腾工业园区"
$ws1.Range("E6").Value = "2024.07.28 09:30-07.28 17:30"
$ws1.Range("F6").Value = 482
$ws1.Range("G6").Value = 60
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86520"
$ws1.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"

# Row 7: 合肥·第二届TH动漫游戏嘉年华
$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = "2024-07-28"
$ws1.Range("C7").Value = "合肥·第二届TH动漫游戏嘉年华"
$ws1.Range("D7").Value = "田埠西路199号 吉祥如意宴会楼蜀山店"
$ws1.Range("E7").Value = "2024.07.28 09:30-07.28 17:00"
$ws1.Range("F7").Value = 86
$ws1.Range("G7").Value = 55
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=87447"
$ws1.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"

# Row 8: 合肥·首届进击的巨人ONLY漫展
$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = "2024-07-28"
$ws1.Range("C8").Value = "合肥·首届进击的巨人ONLY漫展"
$ws1.Range("D8").Value = "胜利路198号 合肥元一希尔顿酒店"
$ws1.Range("E8").Value = "2024.07.28 09:30-07.28 16:30"
$ws1.Range("F8").Value = 98
$ws1.Range("G8").Value = 98
$ws1.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=88965"
$ws1.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202406/q9ANU7gh1718880973689.jpeg"

# Row 9: 巢湖·元气动漫游戏嘉年华
$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "2024-07-30"
$ws1.Range("C9").Value = "巢湖·元气动漫游戏嘉年华"
$ws1.Range("D9").Value = "团结东路7号 巢湖宾馆"
$ws1.Range("E9").Value = "2024.07.30 10:00-07.30 17:00"
$ws1.Range("F9").Value = 34
$ws1.Range("G9").Value = 45
$ws1.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=88193"
$ws1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202406/3VBeQfqQ1719318873395.jpeg"

# Row 10: 合肥·第七届环形宇宙动漫游戏嘉年华
$ws1.Range("A10").Value = 9
$ws1.Range("B10").Value = "2024-08-03"
$ws1.Range("C10").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$ws1.Range("D10").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E10").Value = "2024.08.03 09:30-08.04 17:00"
$ws1.Range("F10").Value = 6335
$ws1.Range("G10").Value = 44.1
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=84767"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"

# Row 11: 合肥·排球少年only之夏日招新季
$ws1.Range("A11").Value = 10
$ws1.Range("B11").Value = "2024-08-10"
$ws1.Range("C11").Value = "合肥·排球少年only之夏日招新季"
$ws1.Range("D11").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws1.Range("E11").Value = "2024.08.10 10:00-08.10 17:00"
$ws1.Range("F11").Value = 213
$ws1.Range("G11").Value = 70
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws1.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"

# Row 12: 合肥·比翼连枝国乙&代号鸢only
$ws1.Range("A12").Value = 11
$ws1.Range("B12").Value = "2024-08-10"
$ws1.Range("C12").Value = "合肥·比翼连枝国乙&代号鸢only"
$ws1.Range("D12").Value = "长江东大街与东二环路交叉口向南300米东方摩域商业广场三楼 格律诗婚礼艺术中心(筑梦店)"
$ws1.Range("E12").Value = "2024.08.10 09:00-08.10 22:00"
$ws1.Range("F12").Value = 337
$ws1.Range("G12").Value = 65
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=88421"
$ws1.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"

# Row 13: 合肥·第八届环形宇宙动漫游戏嘉年华Plus
$ws1.Range("A13").Value = 12
$ws1.Range("B13").Value = "2024-08-17"
$ws1.Range("C13").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws1.Range("D13").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E13").Value = "2024.08.17 09:30-08.18 17:00"
$ws1.Range("F13").Value = 2494
$ws1.Range("G13").Value = 69
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws1.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

# Row 14: 合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票
$ws1.Range("A14").Value = 13
$ws1.Range("B14").Value = "2024-08-17"
$ws1.Range("C14").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws1.Range("D14").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E14").Value = "2024.08.17 09:30-08.17 17:00"
$ws1.Range("F14").Value = 137
$ws1.Range("G14").Value = 0.1
$ws1.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws1.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

# Row 15: 合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票
$ws1.Range("A15").Value = 14
$ws1.Range("B15").Value = "2024-08-17"
$ws1.Range("C15").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws1.Range("D15").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E15").Value = "2024.08.17 09:30-08.17 17:00"
$ws1.Range("F15").Value = 238
$ws1.Range("G15").Value = 0.1
$ws1.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws1.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

# Row 16: 合肥·银魂主题派对only2.0
$ws1.Range("A16").Value = 15
$ws1.Range("B16").Value = "2024-08-17"
$ws1.Range("C16").Value = "合肥·银魂主题派对only2.0"
$ws1.Range("D16").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws1.Range("E16").Value = "2024.08.17 13:00-08.17 18:00"
$ws1.Range("F16").Value = 250
$ws1.Range("G16").Value = 128
$ws1.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws1.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

# Row 17: 合肥·SSS第五人格only
$ws1.Range("A17").Value = 16
$ws1.Range("B17").Value = "2024-08-18"
$ws1.Range("C17").Value = "合肥·SSS第五人格only"
$ws1.Range("D17").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws1.Range("E17").Value = "2024.08.18 09:00-08.18 17:00"
$ws1.Range("F17").Value = 487
$ws1.Range("G17").Value = 68
$ws1.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws1.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"

# --- Sheet 4 ("全部类型") ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows("2:4").Delete()

# Row 2: 合肥·Yolo Fes永乐庆典Vol.3·少女偶像联合演出D
$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "2024-07-26"
$ws4.Range("C2").Value = "合肥·Yolo Fes永乐庆典Vol.3·少女偶像联合演出DAY1&DAY3"
$ws4.Range("D2").Value = "金寨路与天堂窄路交叉口 梵木艺术中心"
$ws4.Range("E2").Value = "2024.07.26 18:00-07.28 23:59"
$ws4.Range("F2").Value = 9
$ws4.Range("G2").Value = 128
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=89514"
$ws4.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202407/aMtLMGR31721289854139.jpeg"

# Row 3: 合肥·Yolo Fes永乐庆典VOL.3·少女偶像联合演出
$ws4.Range("A3").Value = 2
$ws4.Range("B3").Value = "2024-07-27"
$ws4.Range("C3").Value = "合肥·Yolo Fes永乐庆典VOL.3·少女偶像联合演出"
$ws4.Range("D3").Value = "金寨路310号合柴1972文创园区C-1号 合肥九莱福"
$ws4.Range("E3").Value = "2024.07.27 12:00-07.27 23:59"
$ws4.Range("F3").Value = 9
$ws4.Range("G3").Value = 168
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=89435"
$ws4.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202407/gaKvbqHM1720422980894.png"

# Row 4: 合肥·灵能百分百ONLY2.0
$ws4.Range("A4").Value = 3
$ws4.Range("B4").Value = "2024-07-27"
$ws4.Range("C4").Value = "合肥·灵能百分百ONLY2.0"
$ws4.Range("D4").Value = "铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)"
$ws4.Range("E4").Value = "2024.07.27 10:00-07.27 17:00"
$ws4.Range("F4").Value = 117
$ws4.Range("G4").Value = "不可售"
$ws4.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=87497"
$ws4.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg"

# Row 5: 安徽·MAX特摄only展
$ws4.Range("A5").Value = 4
$ws4.Range("B5").Value = "2024-07-27"
$ws4.Range("C5").Value = "安徽·MAX特摄only展"
$ws4.Range("D5").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws4.Range("E5").Value = "2024.07.27 09:30-07.27 18:00"
$ws4.Range("F5").Value = 502
$ws4.Range("G5").Value = 50
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=83684"
$ws4.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"

# Row 6: 庐江·夏日游嘉年华
$ws4.Range("A6").Value = 5
$ws4.Range("B6").Value = "2024-07-27"
$ws4.Range("C6").Value = "庐江·夏日游嘉年华"
$ws4.Range("D6").Value = "白山路东150米 庐江体育馆"
$ws4.Range("E6").Value = "2024.07.27 09:00-07.28 17:00"
$ws4.Range("F6").Value = 183
$ws4.Range("G6").Value = 60
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87569"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg"

# Row 7: 长丰·莓可可游戏动漫展
$ws4.Range("A7").Value = 6
$ws4.Range("B7").Value = "2024-07-27"
$ws4.Range("C7").Value = "长丰·莓可可游戏动漫展"
$ws4.Range("D7").Value = "长寿路12号 长丰宾馆·梅山饭店(长寿路店)"
$ws4.Range("E7").Value = "2024.07.27 10:00-07.27 17:00"
$ws4.Range("F7").Value = 19
$ws4.Range("G7").Value = "不可售"
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=87796"
$ws4.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202406/MLTfeikq1718823574810.png"

# Row 8: 合肥·咒术回战only
$ws4.Range("A8").Value = 7
$ws4.Range("B8").Value = "2024-07-28"
$ws4.Range("C8").Value = "合肥·咒术回战only"
$ws4.Range("D8").Value = "清河路19号 依立腾工业园区"
$ws4.Range("E8").Value = "2024.07.28 09:30-07.28 17:30"
$ws4.Range("F8").Value = 482
$ws4.Range("G8").Value = 60
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=86520"
$ws4.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"

# Row 9: 合肥·第二届TH动漫游戏嘉年华
$ws4.Range("A9").Value = 8
$ws4.Range("B9").Value = "2024-07-28"
$ws4.Range("C9").Value = "合肥·第二届TH动漫游戏嘉年华"
$ws4.Range("D9").Value = "田埠西路199号 吉祥如意宴会楼蜀山店"
$ws4.Range("E9").Value = "2024.07.28 09:30-07.28 17:00"
$ws4.Range("F9").Value = 86
$ws4.Range("G9").Value = 55
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=87447"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"

# Row 10: 合肥·首届进击的巨人ONLY漫展
$ws4.Range("A10").Value = 9
$ws4.Range("B10").Value = "2024-07-28"
$ws4.Range("C10").Value = "合肥·首届进击的巨人ONLY漫展"
$ws4.Range("D10").Value = "胜利路198号 合肥元一希尔顿酒店"
$ws4.Range("E10").Value = "2024.07.28 09:30-07.28 16:30"
$ws4.Range("F10").Value = 98
$ws4.Range("G10").Value = 98
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=88965"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202406/q9ANU7gh1718880973689.jpeg"

# Row 11: 巢湖·元气动漫游戏嘉年华
$ws4.Range("A11").Value = 10
$ws4.Range("B11").Value = "2024-07-30"
$ws4.Range("C11").Value = "巢湖·元气动漫游戏嘉年华"
$ws4.Range("D11").Value = "团结东路7号 巢湖宾馆"
$ws4.Range("E11").Value = "2024.07.30 10:00-07.30 17:00"
$ws4.Range("F11").Value = 34
$ws4.Range("G11").Value = 45
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=88193"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202406/3VBeQfqQ1719318873395.jpeg"

# Row 12: 合肥·新西兰·治愈系民谣歌手Luke Thompson202
$ws4.Range("A12").Value = 11
$ws4.Range("B12").Value = "2024-08-02"
$ws4.Range("C12").Value = "合肥·新西兰·治愈系民谣歌手Luke Thompson2024中国巡演 KEEP ROLLING ON "
$ws4.Range("D12").Value = "宁国路罍街二期15号楼安徽原创音乐基地3楼 合肥ON THE WAY LiveHouse"
$ws4.Range("E12").Value = "2024.08.02 20:00-08.02 21:30"
$ws4.Range("F12").Value = 4
$ws4.Range("G12").Value = 180
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=88824"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202407/FKItJRNl1719803666645.jpeg"

# Row 13: 合肥·第七届环形宇宙动漫游戏嘉年华
$ws4.Range("A13").Value = 12
$ws4.Range("B13").Value = "2024-08-03"
$ws4.Range("C13").Value = "合肥·第七届环形宇宙动漫游戏嘉年华"
$ws4.Range("D13").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E13").Value = "2024.08.03 09:30-08.04 17:00"
$ws4.Range("F13").Value = 6335
$ws4.Range("G13").Value = 44.1
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=84767"
$ws4.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"

# Row 14: 合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻
$ws4.Range("A14").Value = 13
$ws4.Range("B14").Value = "2024-08-03"
$ws4.Range("C14").Value = "合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会"
$ws4.Range("D14").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws4.Range("E14").Value = "2024.08.03 19:30-08.03 21:00"
$ws4.Range("F14").Value = 47
$ws4.Range("G14").Value = 80
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=83556"
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg"

# Row 15: 合肥·排球少年only之夏日招新季
$ws4.Range("A15").Value = 14
$ws4.Range("B15").Value = "2024-08-10"
$ws4.Range("C15").Value = "合肥·排球少年only之夏日招新季"
$ws4.Range("D15").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws4.Range("E15").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("F15").Value = 213
$ws4.Range("G15").Value = 70
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws4.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"

# Row 16: 合肥·比翼连枝国乙&代号鸢only
$ws4.Range("A16").Value = 15
$ws4.Range("B16").Value = "2024-08-10"
$ws4.Range("C16").Value = "合肥·比翼连枝国乙&代号鸢only"
$ws4.Range("D16").Value = "长江东大街与东二环路交叉口向南300米东方摩域商业广场三楼 格律诗婚礼艺术中心(筑梦店)"
$ws4.Range("E16").Value = "2024.08.10 09:00-08.10 22:00"
$ws4.Range("F16").Value = 337
$ws4.Range("G16").Value = 65
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=88421"
$ws4.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202407/RHiXT98J1721199172046.jpeg"

# Row 17: 合肥·第八届环形宇宙动漫游戏嘉年华Plus
$ws4.Range("A17").Value = 16
$ws4.Range("B17").Value = "2024-08-17"
$ws4.Range("C17").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
$ws4.Range("D17").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E17").Value = "2024.08.17 09:30-08.18 17:00"
$ws4.Range("F17").Value = 2494
$ws4.Range("G17").Value = 69
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=88650"
$ws4.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

# Row 18: 合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票
$ws4.Range("A18").Value = 17
$ws4.Range("B18").Value = "2024-08-17"
$ws4.Range("C18").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞周边预约票"
$ws4.Range("D18").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E18").Value = "2024.08.17 09:30-08.17 17:00"
$ws4.Range("F18").Value = 137
$ws4.Range("G18").Value = 0.1
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=89420"
$ws4.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202407/hsiXAged1721203655434.jpeg"

# Row 19: 合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票
$ws4.Range("A19").Value = 18
$ws4.Range("B19").Value = "2024-08-17"
$ws4.Range("C19").Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus~水千丞签售预约票"
$ws4.Range("D19").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E19").Value = "2024.08.17 09:30-08.17 17:00"
$ws4.Range("F19").Value = 238
$ws4.Range("G19").Value = 0.1
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=89421"
$ws4.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202407/r8wJqvVO1721202573195.jpeg"

# Row 20: 合肥·银魂主题派对only2.0
$ws4.Range("A20").Value = 19
$ws4.Range("B20").Value = "2024-08-17"
$ws4.Range("C20").Value = "合肥·银魂主题派对only2.0"
$ws4.Range("D20").Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
$ws4.Range("E20").Value = "2024.08.17 13:00-08.17 18:00"
$ws4.Range("F20").Value = 250
$ws4.Range("G20").Value = 128
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=87173"
$ws4.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

# Row 21: 合肥·SSS第五人格only
$ws4.Range("A21").Value = 20
$ws4.Range("B21").Value = "2024-08-18"
$ws4.Range("C21").Value = "合肥·SSS第五人格only"
$ws4.Range("D21").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws4.Range("E21").Value = "2024.08.18 09:00-08.18 17:00"
$ws4.Range("F21").Value = 487
$ws4.Range("G21").Value = 68
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=88430"
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202406/a0qh8I1h1719660853555.png"
